$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, shifting rows 13-24 down to 14-25.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new weekly price entry.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 4).Value = 44629
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112022
$ws.Cells.Item(13, 7).Value = "Arveja Verde"
$ws.Cells.Item(13, 8).Value = "Perfection"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 35
$ws.Cells.Item(13, 11).Value = 25000
$ws.Cells.Item(13, 12).Value = 26000
$ws.Cells.Item(13, 13).Value = 25429
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 1017
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
